# "food commodity to food group table V4.xlsx" - facet-map update.
#
# The food_groups lookup table (F2:G15) had its "beverages, coffee/tea" key
# (row 3) renamed to "beverages, coffee/tea/cocoa" so that Cocoa (row 24,
# previously classified under "nuts and seeds") can be re-pointed at the
# same bucket as Coffee/Tea ("beverages"). Because the D-column VLOOKUPs use
# an approximate (range) match against this alphabetically sorted table,
# renaming the F3 key also shifts Coffee (row 21) and Tea (row 56) - whose
# food_group_assignment text ("beverages, coffee/tea") no longer exactly
# matches any key - down into the "beverages, alcoholic" bucket, so their
# food_group_code becomes "alcohol". All of this recalculates automatically;
# we only need to edit the two source cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the lookup-table key used by the VLOOKUPs in column D.
$ws.Range("F3").Value = "beverages, coffee/tea/cocoa"

# 2) Re-classify Cocoa's food_group_assignment to match the renamed key.
$ws.Range("C24").Value = "beverages, coffee/tea/cocoa"

# 3) These two rows wrap across two lines now that the text is longer, so
#    Excel auto-sizes them (no explicit user resize -> no customHeight).
$ws.Range("F3").EntireRow.RowHeight = 32
$ws.Range("C24").EntireRow.RowHeight = 32

# 4) Restore the view state: selection moved to I20 and the sheet is
#    scrolled back to the top (topLeftCell reset), matching the saved file.
$ws.Range("I20").Select()
